# bi_oncho_prestop_1_202307_site.xlsx : "bi: update oncho form"
#
# Changes applied:
#  1. choices sheet: the "MPANDA" commune entry is replaced by a new
#     "RUGAZI" commune (still under MPANDA district); the commune list
#     (rows 32-34) is re-sorted alphabetically, which bumps MUGINA and
#     MURWI up a row and puts RUGAZI in their place.
#  2. choices sheet: the centre_sante "KIRENGANE" now belongs to the
#     "RUGAZI" commune (used to be "MPANDA").
#  3. choices sheet: village "MURAMA –RUNGWE" renamed/corrected to
#     "MURAMA RUGWE".
#  4. settings sheet: form_title/form_id bumped from V2 to V3.
#  5. Sheet1 (scratch/staging lookup sheet) mirrors the same data
#     updates in its various helper columns.
#  6. Assorted view-state tweaks (frozen pane, selection, column
#     widths, sort-state range) that Excel recorded as part of the
#     edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1+2+3. choices sheet data updates
# ---------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# Commune list (rows 25-37) - MPANDA commune removed, RUGAZI commune
# inserted, keeping the list alphabetically sorted: MUGINA, MURWI,
# RUGAZI now occupy rows 32-34 (instead of MPANDA, MUGINA, MURWI).
$choices.Range("B32").Value = "MUGINA"
$choices.Range("C32").Value = "MUGINA"
$choices.Range("E32").Value = "MABAYI"

$choices.Range("B33").Value = "MURWI"
$choices.Range("C33").Value = "MURWI"
$choices.Range("E33").Value = "BUKINANYANA"

$choices.Range("B34").Value = "RUGAZI"
$choices.Range("C34").Value = "RUGAZI"
$choices.Range("E34").Value = "MPANDA"

# centre_sante list: KIRENGANE (row 43) now filed under the RUGAZI commune.
$choices.Range("F43").Value = "RUGAZI"

# village_list: spelling/name correction.
$choices.Range("B61").Value = "MURAMA RUGWE"
$choices.Range("C61").Value = "MURAMA RUGWE"

# ---------------------------------------------------------------
# 4. settings sheet: bump form version V2 -> V3
# ---------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "(2023 Juillet) ONCHO pre Stop - 1. Formulaire Site V3"
$settings.Range("B2").Value = "bi_oncho_prestop_1_202307_site_v3"

# ---------------------------------------------------------------
# 5. Sheet1 scratch/staging sheet - mirror the same updates
# ---------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet1")

$sheet1.Range("P6").Value = "RUGAZI"

$sheet1.Range("M9").Value = "MABAYI"
$sheet1.Range("N9").Value = "MUGINA"
$sheet1.Range("T9").Value = "MURAMA RUGWE"

$sheet1.Range("M10").Value = "BUKINANYANA"
$sheet1.Range("N10").Value = "MURWI"

$sheet1.Range("C11").Value = "RUGAZI"
$sheet1.Range("M11").Value = "MPANDA"
$sheet1.Range("N11").Value = "RUGAZI"

$sheet1.Range("E12").Value = "MURAMA RUGWE"

# ---------------------------------------------------------------
# 6. View-state tweaks
# ---------------------------------------------------------------

# choices sheet: frozen pane top-left cell and current selection.
$choices.Range("A44").Select()
$choices.Application.ActiveWindow.FreezePanes = $true
$choices.Range("G54:G68").Select()

# settings sheet: selection.
$settings.Range("B2").Select()

# Sheet1: scrolled view (top-left cell D1) and current selection.
$sheet1.Range("S2:S16").Select()
$sheet1.Application.ActiveWindow.ScrollColumn = 4

# Sheet1: column widths (auto-fit effect of the "MURAMA RUGWE" text change).
$sheet1.Columns.Item(10).ColumnWidth = 16.75
$sheet1.Columns.Item(11).ColumnWidth = 16

# Sheet1: sort-state range shrunk to match the actual data range.
$sheet1.Sort.SortFields.Clear()
$sheet1.Sort.SortFields.Add($sheet1.Range("T2:T16")) | Out-Null
$sheet1.Sort.SetRange($sheet1.Range("S2:T16"))
$sheet1.Sort.Header = -4142
$sheet1.Sort.Apply()
